$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車 (car) sheet

# --- Row 1 (header row): extend with new header labels in H1:N1, and
#     correct C1 label from the old duplicated-data artifact to "capacity" ---

# Style template for header row is style of B1 (bold / bordered / centered)
$headerStyleSrc = $ws.Cells.Item(1, 2)
foreach ($col in 8..14) {
    $dst = $ws.Cells.Item(1, $col)
    $headerStyleSrc.Copy($dst)
}

$ws.Cells.Item(1, 2).Value  = "name"
$ws.Cells.Item(1, 3).Value  = "capacity"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "register_date"
$ws.Cells.Item(1, 6).Value  = "register_reason"
$ws.Cells.Item(1, 7).Value  = "acquire_value"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 & 3 (data rows): extend with new trailing metadata columns
#     H:N, and fix the shifted name/reason values ---

$dataStyleSrc2 = $ws.Cells.Item(2, 2)
foreach ($col in 8..14) {
    $dst = $ws.Cells.Item(2, $col)
    $dataStyleSrc2.Copy($dst)
}

$dataStyleSrc3 = $ws.Cells.Item(3, 2)
foreach ($col in 8..14) {
    $dst = $ws.Cells.Item(3, $col)
    $dataStyleSrc3.Copy($dst)
}

# Row 2
$ws.Cells.Item(2, 2).Value  = "福特六和"
$ws.Cells.Item(2, 3).Value  = 1999
$ws.Cells.Item(2, 4).Value  = "黃偉哲"
$ws.Cells.Item(2, 5).Value  = "98年01月01日"
$ws.Cells.Item(2, 6).Value  = "買賣"
$ws.Cells.Item(2, 7).Value  = 14000
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"
$ws.Cells.Item(2, 10).Value = "2012-05-01"
$ws.Cells.Item(2, 11).Value = "黃偉哲"
$ws.Cells.Item(2, 12).Value = 1367
$ws.Cells.Item(2, 13).Value = "tmp62651"
$ws.Cells.Item(2, 14).Value = 32

# Row 3
$ws.Cells.Item(3, 2).Value  = "曰產"
$ws.Cells.Item(3, 3).Value  = 1995
$ws.Cells.Item(3, 4).Value  = "黃偉哲"
$ws.Cells.Item(3, 5).Value  = "96年06月11闩"
$ws.Cells.Item(3, 6).Value  = "買賣"
$ws.Cells.Item(3, 7).Value  = 100000
$ws.Cells.Item(3, 8).Value  = "land"
$ws.Cells.Item(3, 9).Value  = "normal"
$ws.Cells.Item(3, 10).Value = "2012-05-01"
$ws.Cells.Item(3, 11).Value = "黃偉哲"
$ws.Cells.Item(3, 12).Value = 1367
$ws.Cells.Item(3, 13).Value = "tmp62651"
$ws.Cells.Item(3, 14).Value = 33
